# excel_writer: also include totals for the balance columns
#
# The "Gesamtergebnis" (overall result) sheet's Total row used to leave the
# two balance columns (Startguthaben / Endsaldo) as "N/A" placeholders.
# They should now carry an actual total like every other column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtergebnis")

# Replace the "N/A" placeholders in the balance columns of the Total row
# (row 3: C = Startguthaben, D = Endsaldo) with real numeric totals.
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# Match the row's autofit height now that it only holds numeric data.
$ws.Rows.Item(3).RowHeight = 13.8

# Leave the workbook focused on the Gesamtergebnis sheet / the cell after
# the edited row, matching how the file was left after the edit.
$ws.Activate()
$ws.Range("C6").Select()
